$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 637
$ws.Range("I2").Value = 500.5
$ws.Range("J2").Value = 819
$ws.Range("K2").Value = 500.5
$ws.Range("L2").Value = 819
$ws.Range("M2").Value = -387.5
$ws.Range("N2").Value = -1045
$ws.Range("H51").Value = 2505.5
$ws.Range("I51").Value = 2505.5
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 2505.5
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -2021.5
$ws.Range("N51").ClearContents()
$ws.Range("H70").Value = 3200
$ws.Range("I70").Value = 3200
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 9600
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -9330
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 3200
$ws.Range("I73").Value = 3200
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 9600
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -8664
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 2750
$ws.Range("I80").Value = 2833.3333
$ws.Range("K80").Value = 8499.999899999999
$ws.Range("M80").Value = -7501.999899999999
$ws.Range("H83").Value = 2750
$ws.Range("I83").Value = 2833.3333
$ws.Range("K83").Value = 25499.9997
$ws.Range("M83").Value = -20507.9997
$ws.Range("H86").Value = 3000
$ws.Range("I86").Value = 3000
$ws.Range("K86").Value = 3000
$ws.Range("M86").Value = -1877
$ws.Range("H89").Value = 3000
$ws.Range("I89").Value = 3000
$ws.Range("K89").Value = 15000
$ws.Range("M89").Value = -9384
$ws.Range("H96").Value = 1004.25
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1004.25
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 3012.75
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -5758.75
$ws.Range("H101").Value = 1168
$ws.Range("I101").Value = 1168
$ws.Range("K101").Value = 3504
$ws.Range("M101").Value = -1882

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1155.25
$ws.Range("I2").Value = 1155.25
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1155.25
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1042.25
$ws.Range("N2").ClearContents()
$ws.Range("H40").Value = 12031
$ws.Range("J40").Value = 12031
$ws.Range("L40").Value = 12031
$ws.Range("N40").Value = -12383
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H116").Value = 1155.25
$ws.Range("I116").Value = 1155.25
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1155.25
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1138.75
$ws.Range("N116").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1155.25
$ws.Range("I3").Value = 1155.25
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1155.25
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1041.25
$ws.Range("N3").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 210.33333
$ws.Range("J12").Value = 289.3
$ws.Range("L12").Value = 867.9000000000001
$ws.Range("N12").Value = -1213.9
$ws.Range("H15").Value = 121.3
$ws.Range("I15").Value = 104.28571
$ws.Range("J15").Value = 161
$ws.Range("K15").Value = 312.85713
$ws.Range("L15").Value = 483
$ws.Range("M15").Value = -172.85713
$ws.Range("N15").Value = -763
$ws.Range("H107").Value = 242.8
$ws.Range("I107").Value = 234.5
$ws.Range("J107").Value = 248.33333
$ws.Range("K107").Value = 703.5
$ws.Range("L107").Value = 744.99999
$ws.Range("M107").Value = 1216.5
$ws.Range("N107").Value = -4584.99999
$ws.Range("H120").Value = 7500
$ws.Range("I120").Value = 3333.3333
$ws.Range("K120").Value = 9999.999899999999
$ws.Range("M120").Value = -5161.999899999999
$ws.Range("H138").Value = 303399.6
$ws.Range("I138").Value = 303399.6
$ws.Range("K138").Value = 910198.7999999999
$ws.Range("M138").Value = -905058.7999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()
$ws.Range("H40").Value = 510999.7
$ws.Range("I40").Value = 14799.6
$ws.Range("J40").Value = 1007199.8
$ws.Range("K40").Value = 14799.6
$ws.Range("L40").Value = 1007199.8
$ws.Range("M40").Value = -14663.6
$ws.Range("N40").Value = -1007471.8
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H61").Value = 4499.5
$ws.Range("I61").Value = 4499.5
$ws.Range("K61").Value = 4499.5
$ws.Range("M61").Value = -4297.5
$ws.Range("H93").Value = 1796
$ws.Range("I93").Value = 1796
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1796
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -548
$ws.Range("N93").ClearContents()
$ws.Range("H113").Value = 4499.5
$ws.Range("I113").Value = 4499.5
$ws.Range("K113").Value = 4499.5
$ws.Range("M113").Value = -2329.5
$ws.Range("H135").Value = 199999
$ws.Range("J135").Value = 199999
$ws.Range("L135").Value = 199999
$ws.Range("N135").Value = -210139
$ws.Range("H136").Value = 1000730
$ws.Range("I136").Value = 1000730
$ws.Range("K136").Value = 3002190
$ws.Range("M136").Value = -2999640

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 50263
$ws.Range("I64").Value = 90000
$ws.Range("K64").Value = 90000
$ws.Range("M64").Value = -89752
$ws.Range("H67").Value = 50263
$ws.Range("I67").Value = 90000
$ws.Range("K67").Value = 90000
$ws.Range("M67").Value = -89142
$ws.Range("H100").Value = 625.2857
$ws.Range("I100").Value = 525.4
$ws.Range("J100").Value = 875
$ws.Range("K100").Value = 1050.8
$ws.Range("L100").Value = 1750
$ws.Range("M100").Value = -509.8
$ws.Range("N100").Value = -2832
$ws.Range("H122").Value = 994
$ws.Range("J122").Value = 993.5
$ws.Range("L122").Value = 2980.5
$ws.Range("N122").Value = -7880.5
$ws.Range("H136").Value = 1976.1333
$ws.Range("I136").Value = 1822.091
$ws.Range("K136").Value = 5466.272999999999
$ws.Range("M136").Value = -2916.272999999999
